$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# New data rows: Hill climber with small (row 8) / Hill climber with big (row 9)
# ---------------------------------------------------------------------------
$ws.Range("C8").Value = 227
$ws.Range("D8").Value = 15
$ws.Range("E8").Value = 29
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 53
$ws.Range("H8").Value = 79
$ws.Range("I8").Value = 19
$ws.Range("J8").Value = 51
$ws.Range("K8").Value = 10
$ws.Range("L8").Value = 34
$ws.Range("N8").Formula = "=MIN(C8:L8)"
$ws.Range("O8").Formula = "=MAX(C8:L8)"
$ws.Range("P8").Formula = "=AVERAGE(C8:L8)"

$ws.Range("C9").Value = 22
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = 64
$ws.Range("F9").Value = 76
$ws.Range("G9").Value = 104
$ws.Range("H9").Value = 45
$ws.Range("I9").Value = 22
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 63
$ws.Range("L9").Value = 144
$ws.Range("N9").Formula = "=MIN(C9:L9)"
$ws.Range("O9").Formula = "=MAX(C9:L9)"
$ws.Range("P9").Formula = "=AVERAGE(C9:L9)"

# ---------------------------------------------------------------------------
# New "outperform" comparison columns (R label / S percentage)
# ---------------------------------------------------------------------------
$ws.Range("R2").Value = "GA outperform random with small"
$ws.Range("S2").Value = 0.8

$ws.Range("R3").Value = "GA outperform hill climber with small"
$ws.Range("S3").Value = 0.4

$ws.Range("R5").Value = "GA outperform random with big"
$ws.Range("S5").Value = 0.8

$ws.Range("R6").Value = "GA outperform hill climber with big"
$ws.Range("S6").Value = 0.6

$ws.Range("R8").Value = "Hill climber outperform random with small"
$ws.Range("S8").Value = 0.6

$ws.Range("R9").Value = "Hill climber outperform random with big"
$ws.Range("S9").Value = 0.6

# ---------------------------------------------------------------------------
# Styling: "Heading 1" cell style for header row + row-label cells
# ---------------------------------------------------------------------------
$ws.Range("C1:L1").Style = "Heading 1"
$ws.Range("N1:P1").Style = "Heading 1"
$ws.Range("A2").Style = "Heading 1"
$ws.Range("A3").Style = "Heading 1"
$ws.Range("A5").Style = "Heading 1"
$ws.Range("A6").Style = "Heading 1"
$ws.Range("A8").Style = "Heading 1"
$ws.Range("A9").Style = "Heading 1"

# Percentage formatting for the new "outperform" ratio column
$ws.Range("S2").NumberFormat = "0%"
$ws.Range("S3").NumberFormat = "0%"
$ws.Range("S5").NumberFormat = "0%"
$ws.Range("S6").NumberFormat = "0%"
$ws.Range("S8").NumberFormat = "0%"
$ws.Range("S9").NumberFormat = "0%"

# ---------------------------------------------------------------------------
# Row heights (header/band rows taller, spacer rows shorter)
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 20.25
$ws.Rows.Item(2).RowHeight = 21
$ws.Rows.Item(3).RowHeight = 21
$ws.Range("A4:S4").EntireRow.RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 20.25
$ws.Rows.Item(6).RowHeight = 21
$ws.Range("A7:S7").EntireRow.RowHeight = 15.75
$ws.Rows.Item(8).RowHeight = 20.25
$ws.Rows.Item(9).RowHeight = 21
$ws.Range("A10:S10").EntireRow.RowHeight = 15.75

# ---------------------------------------------------------------------------
# Column widths (approximate best-fit sizing for the new/expanded columns)
# ---------------------------------------------------------------------------
$ws.Range("A1").ColumnWidth = 27.7
$ws.Range("C1").ColumnWidth = 7.3
$ws.Range("D1").ColumnWidth = 7.3
$ws.Range("E1").ColumnWidth = 7.3
$ws.Range("F1").ColumnWidth = 7.3
$ws.Range("G1").ColumnWidth = 7.3
$ws.Range("H1").ColumnWidth = 7.3
$ws.Range("I1").ColumnWidth = 7.3
$ws.Range("J1").ColumnWidth = 7.3
$ws.Range("K1").ColumnWidth = 7.3
$ws.Range("L1").ColumnWidth = 8.7
$ws.Range("N1").ColumnWidth = 8.9
$ws.Range("O1").ColumnWidth = 9.45
$ws.Range("P1").ColumnWidth = 10
$ws.Range("R1").ColumnWidth = 38.9

# ---------------------------------------------------------------------------
# Selection (matches the saved workbook view)
# ---------------------------------------------------------------------------
[void]$ws.Range("K18").Select()
